$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 12.73660532795967
$ws.Range("B2").Value = 0.0000000000000001110223024625157
$ws.Range("C2").Value = 0.008330174690077187
$ws.Range("D2").Value = 1.225150170186294
$ws.Range("E2").Value = 1.500992939507504

$ws.Range("A3").Value = 9.604651620884498
$ws.Range("B3").Value = 0.0000000000000001110223024625157
$ws.Range("C3").Value = 0.006940735173005978
$ws.Range("D3").Value = 1.020800066600701
$ws.Range("E3").Value = 1.042032775971996

$ws.Range("A4").Value = 11.96895982805127
$ws.Range("B4").Value = 0.0000000000000001110223024625157
$ws.Range("C4").Value = 0.008223955695972059
$ws.Range("D4").Value = 1.209528142612259
$ws.Range("E4").Value = 1.46295832777106

$ws.Range("A5").Value = 11.08076194571725
$ws.Range("B5").Value = 0.0000000000000001110223024625157
$ws.Range("C5").Value = 0.008097047406838982
$ws.Range("D5").Value = 1.190863262484997
$ws.Range("E5").Value = 1.41815530993641

$ws.Range("A6").Value = 8.947102642240461
$ws.Range("B6").Value = 0.0000000000000001110223024625157
$ws.Range("C6").Value = 0.006572778634904253
$ws.Range("D6").Value = 0.9666833125050783
$ws.Range("E6").Value = 0.934476626675791

